$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.265.08"
$ws.Range("E2").Value = "  -3.14%  "
$ws.Range("D3").Value = "2.289.12"
$ws.Range("E3").Value = "  -3.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "494.44"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.24"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("D9").Value = "2.290.67"
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0939"
$ws.Range("E10").Value = "  -3.67%  "
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.75"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").Value = "2.695.03"
$ws.Range("E14").Value = "  -3.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.40"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "54.243.54"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "2.274.27"
$ws.Range("E18").Value = "  -6.30%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.02"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.69"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "303.73"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.17"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.32"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("D30").Value = "0.0₃0701"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.80"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.995"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.08"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.63"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.851"
$ws.Range("E38").Value = "  +6.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.64"
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.78"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.39"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.57"
$ws.Range("E44").Value = "  -5.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.73"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0884"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "240.09"
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0205"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.53"
$ws.Range("E51").Value = "  -2.04%  "
